# Applies the edit described by the diff:
#   1. Drop the "_GoBack" bookmark that used to sit right after
#      "arXiv:1811.04422v1".
#   2. In the "Top result is not the query in updated results." paragraph
#      that follows "arXiv:1811.04422v1", insert the missing word "the"
#      (splitting the single run into three runs).
#   3. In the "Top result is not the query in updated results." paragraph
#      that follows "arXiv:1708.06939", insert "the " the same way, and
#      re-create the "_GoBack" bookmark between the "the " run and the
#      "updated results." run.
#   The third occurrence of that sentence (after "arXiv:1811.11402v2")
#   is left untouched.

$d = $word.ActiveDocument

$targetText = "Top result is not the query in updated results."

# Locate the (1-based) paragraph indices of interest by matching text,
# rather than hard-coding numbers. Character offsets are *not* captured
# here because later edits shift them; only paragraph indices are kept
# (paragraph count never changes in this script, so indices stay valid).
$count = $d.Paragraphs.Count

$paraIndexNoBookmark = -1
$paraIndexWithBookmark = -1
$prevText = ""

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.Trim()
    if ($t -eq $targetText) {
        if ($prevText -eq "arXiv:1811.04422v1") {
            $paraIndexNoBookmark = $i
        } elseif ($prevText -like "arXiv:1708.06939*") {
            $paraIndexWithBookmark = $i
        }
    }
    $prevText = $t
}

# --- 1. Remove the existing "_GoBack" bookmark -----------------------
# (collapsed bookmark, so this never shifts character offsets)
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# Process the *later* paragraph in the document first, then the earlier
# one, so that inserting text for one target never shifts the absolute
# character offsets still needed for the other.

# --- 2. Second occurrence (after "arXiv:1708.06939"): "the " inserted
#        (trailing space), and the "_GoBack" bookmark is re-created
#        right after it.
#        "Top result is not the query in " | "the " | "updated results."
if ($paraIndexWithBookmark -gt 0) {
    $startWithBookmark = $d.Paragraphs($paraIndexWithBookmark).Range.Start
    $insertionOffset2 = 31
    $ins2 = $d.Range($startWithBookmark + $insertionOffset2, $startWithBookmark + $insertionOffset2)
    $ins2.InsertAfter("the ")

    $newRun2 = $d.Range($startWithBookmark + $insertionOffset2, $startWithBookmark + $insertionOffset2 + 4)
    $newRun2.Bold = 1
    $newRun2.Bold = 0

    $bookmarkPoint = $startWithBookmark + $insertionOffset2 + 4
    $bookmarkRange = $d.Range($bookmarkPoint, $bookmarkPoint)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
}

# --- 3. First occurrence (after "arXiv:1811.04422v1"): no bookmark,
#        "the" inserted without extra surrounding spaces (" the" is the
#        new run, no trailing space).
#        "Top result is not the query in" | " the" | " updated results."
if ($paraIndexNoBookmark -gt 0) {
    $startNoBookmark = $d.Paragraphs($paraIndexNoBookmark).Range.Start
    $insertionOffset = 30
    $ins = $d.Range($startNoBookmark + $insertionOffset, $startNoBookmark + $insertionOffset)
    $ins.InsertAfter(" the")

    $newRun = $d.Range($startNoBookmark + $insertionOffset, $startNoBookmark + $insertionOffset + 4)
    $newRun.Bold = 1
    $newRun.Bold = 0
}
